# Re-added sleeve bearing for tonearm; includes tonearm rod updates
#
# Adds a new BOM line (row 43) for the bronze sleeve bearing used on the
# tonearm rod, widens the "Details" column to fit the new text, and
# restores the viewport/selection state left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Details") needs to be a bit wider to fit the new part's notes.
# This engine's ColumnWidth setter adds a constant ~0.8333 char offset on
# save relative to the raw character-width value, so back it out here to
# land on exactly 50.
$ws.Range("B:B").ColumnWidth = 49.166666666666664

# New row 43: the re-added sleeve bearing part.
$ws.Range("A43").Value = "High-Temperature Dry-Running 841 Bronze Sleeve Bearing"
$ws.Range("B43").Value = "for 3/16`" Shaft Diameter and 1/4`" Housing ID, 1/4`" Long"
$ws.Range("D43").Value = 1
$ws.Range("F43").Value = 0.61
$ws.Range("G43").Value = "McMaster"

# Hyperlink the part name to its McMaster-Carr product page, matching the
# style used by every other "Part Name/URL" entry in column A.
$ws.Hyperlinks.Add($ws.Range("A43"), "https://www.mcmaster.com/6338K130/")
$ws.Range("A43").Style = "Hyperlink"

# Restore the scroll position / active selection left in the sheet after
# the edit.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
